# TC39_Canine_Filter_Breed-Rottweiler.xlsx
# "Fixed variables and query errors in Bread from TC30 to TC47"
#
# The CasesTab Cypher query (column B, row 2 on the "startup" sheet) referenced
# an `OPTIONAL MATCH (co:cohort)...` / `co.cohort_description` pair that isn't
# returned correctly together with the rest of the RETURN clause - the trailing
# ", coalesce(co.cohort_description, '') AS `Cohort`" column is removed so the
# query now ends cleanly on the "Response to Treatment" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$casesTabQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed  IN ['Rottweiler']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $casesTabQuery

# The query text got two lines shorter, so the wrapped-text row shrinks too.
$ws.Rows.Item(2).RowHeight = 244.8

# Leave the cursor on the cell that was just corrected.
$ws.Range("B2").Select() | Out-Null
